$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Magic")
$ws.Activate()

$ws.Range("V6").Value = 1.1
$ws.Range("Y6").Value = 2

$ws.Range("T7").Value = 1
$ws.Range("V7").Value = 35000
$ws.Range("W7").Value = 35000
$ws.Range("X7").Value = 35000
$ws.Range("Y7").Value = 35000
$ws.Range("Z7").Value = 35000
